$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting for new rows 207-213 from existing rows (A column: bold/border/center style; E column: date format)
$ws.Range("A206").Copy() | Out-Null
$ws.Range("A207:A213").PasteSpecial(-4122) | Out-Null
$ws.Range("E206").Copy() | Out-Null
$ws.Range("E207:E213").PasteSpecial(-4122) | Out-Null

# --- Swap rows 2 and 3 (columns B, F:AC); C/D/E unchanged ---
    $ws.Range("B2").Value = 6720844
    $ws.Range("F2").Value = 'Guarani Asuncion'
    $ws.Range("G2").Value = 'Olimpia Asuncion'
    $ws.Range("H2").Value = 1
    $ws.Range("I2").Value = 2
    $ws.Range("J2").Value = 'A'
    $ws.Range("K2").Value = 2.45
    $ws.Range("L2").Value = 3
    $ws.Range("M2").Value = 2.75
    $ws.Range("N2").Value = 4
    $ws.Range("O2").Value = 3.2
    $ws.Range("P2").Value = 1.85
    $ws.Range("Q2").Value = 0.5
    $ws.Range("R2").Value = 1.875
    $ws.Range("S2").Value = 1.925
    $ws.Range("T2").Value = 2.5
    $ws.Range("U2").Value = 1.925
    $ws.Range("V2").Value = 1.875
    $ws.Range("W2").Value = -1
    $ws.Range("X2").Value = -1
    $ws.Range("Y2").Value = 0.8500000000000001
    $ws.Range("Z2").Value = -1
    $ws.Range("AA2").Value = 0.925
    $ws.Range("AB2").Value = 0.925
    $ws.Range("AC2").Value = -1
    $ws.Range("B3").Value = 6720873
    $ws.Range("F3").Value = 'Sportivo Luqueno'
    $ws.Range("G3").Value = 'Sportivo Trinidense'
    $ws.Range("H3").Value = 2
    $ws.Range("I3").Value = 2
    $ws.Range("J3").Value = 'D'
    $ws.Range("K3").Value = 2.625
    $ws.Range("L3").Value = 3.1
    $ws.Range("M3").Value = 2.5
    $ws.Range("N3").Value = 2.3
    $ws.Range("O3").Value = 3.1
    $ws.Range("P3").Value = 2.9
    $ws.Range("Q3").Value = -0.25
    $ws.Range("R3").Value = 2.025
    $ws.Range("S3").Value = 1.775
    $ws.Range("T3").Value = 2.5
    $ws.Range("U3").Value = 1.95
    $ws.Range("V3").Value = 1.85
    $ws.Range("W3").Value = -1
    $ws.Range("X3").Value = 2.1
    $ws.Range("Y3").Value = -1
    $ws.Range("Z3").Value = -0.5
    $ws.Range("AA3").Value = 0.3875
    $ws.Range("AB3").Value = 0.95
    $ws.Range("AC3").Value = -1

# --- Cyclic permutation rows 143 -> 144 -> 145 -> 143 ---
# new143 = old145, new144 = old143, new145 = old144
    $ws.Range("B143").Value = 7493433
    $ws.Range("F143").Value = 'Sportivo Luqueno'
    $ws.Range("G143").Value = 'Nacional Asuncion'
    $ws.Range("H143").Value = 1
    $ws.Range("I143").Value = 1
    $ws.Range("J143").Value = 'D'
    $ws.Range("K143").Value = 2.75
    $ws.Range("L143").Value = 3.2
    $ws.Range("M143").Value = 2.4
    $ws.Range("N143").Value = 2.75
    $ws.Range("O143").Value = 3.1
    $ws.Range("P143").Value = 2.45
    $ws.Range("Q143").Value = 0.25
    $ws.Range("R143").Value = 1.75
    $ws.Range("S143").Value = 2.05
    $ws.Range("T143").Value = 2.25
    $ws.Range("U143").Value = 2
    $ws.Range("V143").Value = 1.8
    $ws.Range("W143").Value = -1
    $ws.Range("X143").Value = 2.1
    $ws.Range("Y143").Value = -1
    $ws.Range("Z143").Value = 0.375
    $ws.Range("AA143").Value = -0.5
    $ws.Range("AB143").Value = -0.5
    $ws.Range("AC143").Value = 0.4
    $ws.Range("B144").Value = 7493312
    $ws.Range("F144").Value = 'Cerro Porteno'
    $ws.Range("G144").Value = 'Guarani Asuncion'
    $ws.Range("H144").Value = 4
    $ws.Range("I144").Value = 0
    $ws.Range("J144").Value = 'H'
    $ws.Range("K144").Value = 1.7
    $ws.Range("L144").Value = 3.6
    $ws.Range("M144").Value = 4.333
    $ws.Range("N144").Value = 1.727
    $ws.Range("O144").Value = 3.75
    $ws.Range("P144").Value = 4.2
    $ws.Range("Q144").Value = -0.5
    $ws.Range("R144").Value = 1.8
    $ws.Range("S144").Value = 2
    $ws.Range("T144").Value = 2.75
    $ws.Range("U144").Value = 1.875
    $ws.Range("V144").Value = 1.925
    $ws.Range("W144").Value = 0.7270000000000001
    $ws.Range("X144").Value = -1
    $ws.Range("Y144").Value = -1
    $ws.Range("Z144").Value = 0.8
    $ws.Range("AA144").Value = -1
    $ws.Range("AB144").Value = 0.875
    $ws.Range("AC144").Value = -1
    $ws.Range("B145").Value = 7493311
    $ws.Range("F145").Value = 'General Caballero JLM'
    $ws.Range("G145").Value = 'Olimpia Asuncion'
    $ws.Range("H145").Value = 0
    $ws.Range("I145").Value = 1
    $ws.Range("J145").Value = 'A'
    $ws.Range("K145").Value = 3.4
    $ws.Range("L145").Value = 3.3
    $ws.Range("M145").Value = 2
    $ws.Range("N145").Value = 3.2
    $ws.Range("O145").Value = 3.25
    $ws.Range("P145").Value = 2.1
    $ws.Range("Q145").Value = 0.25
    $ws.Range("R145").Value = 1.95
    $ws.Range("S145").Value = 1.85
    $ws.Range("T145").Value = 2.25
    $ws.Range("U145").Value = 1.775
    $ws.Range("V145").Value = 2.025
    $ws.Range("W145").Value = -1
    $ws.Range("X145").Value = -1
    $ws.Range("Y145").Value = 1.1
    $ws.Range("Z145").Value = -1
    $ws.Range("AA145").Value = 0.8500000000000001
    $ws.Range("AB145").Value = -1
    $ws.Range("AC145").Value = 1.025
# --- Row 206 update + new rows 207-213 ---
# Row 206
$ws.Range("A206").Value = 204
$ws.Range("B206").Value = 7609201
$ws.Range("C206").Value = 'Paraguay Division Profesional'
$ws.Range("D206").Value = 'Paraguay Division Profesional'
$ws.Range("E206").Value = 45381.79166666666
$ws.Range("F206").Value = 'Cerro Porteno'
$ws.Range("G206").Value = 'Nacional Asuncion'
$ws.Range("H206").Value = 1
$ws.Range("I206").Value = 0
$ws.Range("J206").Value = 'H'
$ws.Range("K206").Value = 1.5
$ws.Range("L206").Value = 4
$ws.Range("M206").Value = 6
$ws.Range("N206").Value = 1.4
$ws.Range("O206").Value = 4.2
$ws.Range("P206").Value = 7
$ws.Range("Q206").Value = -1
$ws.Range("R206").Value = 1.775
$ws.Range("S206").Value = 2.025
$ws.Range("T206").Value = 2.5
$ws.Range("U206").Value = 1.85
$ws.Range("V206").Value = 1.95
$ws.Range("W206").Value = 0.3999999999999999
$ws.Range("X206").Value = -1
$ws.Range("Y206").Value = -1
$ws.Range("Z206").Value = 0
$ws.Range("AA206").Value = 0
$ws.Range("AB206").Value = -1
$ws.Range("AC206").Value = 0.95

# Row 207
$ws.Range("A207").Value = 205
$ws.Range("B207").Value = 7609146
$ws.Range("C207").Value = 'Paraguay Division Profesional'
$ws.Range("D207").Value = 'Paraguay Division Profesional'
$ws.Range("E207").Value = 45381.89583333334
$ws.Range("F207").Value = 'Libertad Asuncion'
$ws.Range("G207").Value = 'Olimpia Asuncion'
$ws.Range("H207").Value = 0
$ws.Range("I207").Value = 0
$ws.Range("J207").Value = 'D'
$ws.Range("K207").Value = 1.8
$ws.Range("L207").Value = 3.3
$ws.Range("M207").Value = 4
$ws.Range("N207").Value = 1.85
$ws.Range("O207").Value = 3.25
$ws.Range("P207").Value = 3.8
$ws.Range("Q207").Value = -0.5
$ws.Range("R207").Value = 1.95
$ws.Range("S207").Value = 1.85
$ws.Range("T207").Value = 2.25
$ws.Range("U207").Value = 1.85
$ws.Range("V207").Value = 1.95
$ws.Range("W207").Value = -1
$ws.Range("X207").Value = 2.25
$ws.Range("Y207").Value = -1
$ws.Range("Z207").Value = -1
$ws.Range("AA207").Value = 0.8500000000000001
$ws.Range("AB207").Value = -1
$ws.Range("AC207").Value = 0.95

# Row 208
$ws.Range("A208").Value = 206
$ws.Range("B208").Value = 7609145
$ws.Range("C208").Value = 'Paraguay Division Profesional'
$ws.Range("D208").Value = 'Paraguay Division Profesional'
$ws.Range("E208").Value = 45382.77083333334
$ws.Range("F208").Value = 'Sol de America'
$ws.Range("G208").Value = 'General Caballero JLM'
$ws.Range("H208").Value = 2
$ws.Range("I208").Value = 1
$ws.Range("J208").Value = 'H'
$ws.Range("K208").Value = 2.1
$ws.Range("L208").Value = 3.2
$ws.Range("M208").Value = 3.2
$ws.Range("N208").Value = 2.25
$ws.Range("O208").Value = 3.2
$ws.Range("P208").Value = 3
$ws.Range("Q208").Value = -0.25
$ws.Range("R208").Value = 1.9
$ws.Range("S208").Value = 1.9
$ws.Range("T208").Value = 2.25
$ws.Range("U208").Value = 1.925
$ws.Range("V208").Value = 1.875
$ws.Range("W208").Value = 1.25
$ws.Range("X208").Value = -1
$ws.Range("Y208").Value = -1
$ws.Range("Z208").Value = 0.8999999999999999
$ws.Range("AA208").Value = -1
$ws.Range("AB208").Value = 0.925
$ws.Range("AC208").Value = -1

# Row 209
$ws.Range("A209").Value = 207
$ws.Range("B209").Value = 7609672
$ws.Range("C209").Value = 'Paraguay Division Profesional'
$ws.Range("D209").Value = 'Paraguay Division Profesional'
$ws.Range("E209").Value = 45382.875
$ws.Range("F209").Value = '2 de Mayo'
$ws.Range("G209").Value = 'Guarani Asuncion'
$ws.Range("H209").Value = 0
$ws.Range("I209").Value = 1
$ws.Range("J209").Value = 'A'
$ws.Range("K209").Value = 2.2
$ws.Range("L209").Value = 3.1
$ws.Range("M209").Value = 3.1
$ws.Range("N209").Value = 2.55
$ws.Range("O209").Value = 3
$ws.Range("P209").Value = 2.7
$ws.Range("Q209").Value = 0
$ws.Range("R209").Value = 1.85
$ws.Range("S209").Value = 1.95
$ws.Range("T209").Value = 2.25
$ws.Range("U209").Value = 2
$ws.Range("V209").Value = 1.8
$ws.Range("W209").Value = -1
$ws.Range("X209").Value = -1
$ws.Range("Y209").Value = 1.7
$ws.Range("Z209").Value = -1
$ws.Range("AA209").Value = 0.95
$ws.Range("AB209").Value = -1
$ws.Range("AC209").Value = 0.8

# Row 210
$ws.Range("A210").Value = 208
$ws.Range("B210").Value = 7609151
$ws.Range("C210").Value = 'Paraguay Division Profesional'
$ws.Range("D210").Value = 'Paraguay Division Profesional'
$ws.Range("E210").Value = 45387.79166666666
$ws.Range("F210").Value = 'Tacuary'
$ws.Range("G210").Value = 'General Caballero JLM'
$ws.Range("K210").Value = 2.55
$ws.Range("L210").Value = 3.2
$ws.Range("M210").Value = 2.7
$ws.Range("N210").Value = 2.55
$ws.Range("O210").Value = 3.2
$ws.Range("P210").Value = 2.7
$ws.Range("Q210").Value = 0
$ws.Range("R210").Value = 1.825
$ws.Range("S210").Value = 1.975
$ws.Range("T210").Value = 2.25
$ws.Range("U210").Value = 1.875
$ws.Range("V210").Value = 1.925
$ws.Range("W210").Value = 0
$ws.Range("X210").Value = 0
$ws.Range("Y210").Value = 0
$ws.Range("Z210").Value = 0
$ws.Range("AA210").Value = 0

# Row 211
$ws.Range("A211").Value = 209
$ws.Range("B211").Value = 7609203
$ws.Range("C211").Value = 'Paraguay Division Profesional'
$ws.Range("D211").Value = 'Paraguay Division Profesional'
$ws.Range("E211").Value = 45387.89583333334
$ws.Range("F211").Value = 'Nacional Asuncion'
$ws.Range("G211").Value = 'Sportivo Ameliano'
$ws.Range("K211").Value = 2.625
$ws.Range("L211").Value = 3.2
$ws.Range("M211").Value = 2.6
$ws.Range("N211").Value = 2.625
$ws.Range("O211").Value = 3.2
$ws.Range("P211").Value = 2.6
$ws.Range("Q211").Value = 0
$ws.Range("R211").Value = 1.95
$ws.Range("S211").Value = 1.85
$ws.Range("T211").Value = 2.25
$ws.Range("U211").Value = 1.9
$ws.Range("V211").Value = 1.9
$ws.Range("W211").Value = 0
$ws.Range("X211").Value = 0
$ws.Range("Y211").Value = 0
$ws.Range("Z211").Value = 0
$ws.Range("AA211").Value = 0

# Row 212
$ws.Range("A212").Value = 210
$ws.Range("B212").Value = 7609148
$ws.Range("C212").Value = 'Paraguay Division Profesional'
$ws.Range("D212").Value = 'Paraguay Division Profesional'
$ws.Range("E212").Value = 45388.79166666666
$ws.Range("F212").Value = 'Sportivo Trinidense'
$ws.Range("G212").Value = 'Cerro Porteno'
$ws.Range("K212").Value = 6
$ws.Range("L212").Value = 4.5
$ws.Range("M212").Value = 1.444
$ws.Range("N212").Value = 6
$ws.Range("O212").Value = 4.5
$ws.Range("P212").Value = 1.444
$ws.Range("Q212").Value = 1.25
$ws.Range("R212").Value = 1.825
$ws.Range("S212").Value = 1.975
$ws.Range("T212").Value = 2.75
$ws.Range("U212").Value = 1.825
$ws.Range("V212").Value = 1.975
$ws.Range("W212").Value = 0
$ws.Range("X212").Value = 0
$ws.Range("Y212").Value = 0
$ws.Range("Z212").Value = 0
$ws.Range("AA212").Value = 0

# Row 213
$ws.Range("A213").Value = 211
$ws.Range("B213").Value = 7609149
$ws.Range("C213").Value = 'Paraguay Division Profesional'
$ws.Range("D213").Value = 'Paraguay Division Profesional'
$ws.Range("E213").Value = 45388.89583333334
$ws.Range("F213").Value = 'Libertad Asuncion'
$ws.Range("G213").Value = 'Guarani Asuncion'
$ws.Range("K213").Value = 1.5
$ws.Range("L213").Value = 4
$ws.Range("M213").Value = 6
$ws.Range("N213").Value = 1.5
$ws.Range("O213").Value = 4
$ws.Range("P213").Value = 6
$ws.Range("Q213").Value = -1
$ws.Range("R213").Value = 1.85
$ws.Range("S213").Value = 1.95
$ws.Range("T213").Value = 2.5
$ws.Range("U213").Value = 1.875
$ws.Range("V213").Value = 1.925
$ws.Range("W213").Value = 0
$ws.Range("X213").Value = 0
$ws.Range("Y213").Value = 0
$ws.Range("Z213").Value = 0
$ws.Range("AA213").Value = 0
